# Change 1: "Used when the dependent variable is binary in nature. (or categorical)"
#   -> "...(or categorical – not quantitative in nature)" with a _GoBack bookmark
#      inserted right before the closing paren, split across three runs.
$d = $word.ActiveDocument

$r = $d.Content
$null = $r.Find.Execute("(or categorical)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$closeParenPos = $r.End - 1
$insertPoint = $d.Range($closeParenPos, $closeParenPos)
$insertPoint.InsertBefore(" – not quantitative in nature")

# Force the newly inserted text onto its own run (same formatting as before/after)
# by toggling a formatting property on and back off - this breaks run-coalescing
# without leaving any residual formatting mark in the saved XML.
$rNew = $d.Content
$null = $rNew.Find.Execute(" – not quantitative in nature", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rNew.Bold = 1
$rNew2 = $d.Content
$null = $rNew2.Find.Execute(" – not quantitative in nature", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rNew2.Bold = 0
$endOfInsertedText = $rNew2.End

# Move the document's _GoBack bookmark here (right before the trailing ")").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($endOfInsertedText, $endOfInsertedText)
$d.Bookmarks.Add("_GoBack", $bmRange)


# Change 2: "The adjust predictions and probabilities should match the distribution
#   of an observed set of labels" -> "The adjusted predictions and probabilities..."
#   (insert "ed" after "adjust", split into three runs, all keeping the existing
#   italic formatting of the source run).
$r2 = $d.Content
$null = $r2.Find.Execute("The adjust predictions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos2 = $r2.Start + 10
$ip2 = $d.Range($insertPos2, $insertPos2)
$ip2.InsertBefore("ed")

$rEd = $d.Content
$null = $rEd.Find.Execute("adjusted predictions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$edStart = $rEd.Start + 6
$edEnd = $edStart + 2
$edRange = $d.Range($edStart, $edEnd)
$edRange.Bold = 1
$edRange2 = $d.Range($edStart, $edEnd)
$edRange2.Bold = 0


# Change 3: the stray _GoBack bookmark that used to sit after the MNIST paragraph
# was already relocated above (Bookmarks collection only allows one bookmark per
# name), so nothing further is required here - it no longer trails the MNIST
# paragraph's run.

Write-Output "edits applied"
